# Chapter04/Normal.xlsx — add the normal-distribution analysis formulas
# (mean, NORM.DIST / NORM.INV lookups, and the per-row density column)
# that drive the scatter chart already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: recompute the mean from the sample data instead of the literal 100
$ws.Range("B1").Formula = "=AVERAGE(A4:A68)"

# E1/E2: probability of exactly 84 and probability of 84-or-less
$ws.Range("E1").Formula = "=NORM.DIST(84, B1, B2, FALSE)"
$ws.Range("E2").Formula = "=NORM.DIST(84, B1, B2, TRUE)"

# H1/H2: value below which 71% / 85% of the distribution falls
$ws.Range("H1").Formula = "=NORM.INV(0.71, B1, B2)"
$ws.Range("H2").Formula = "=NORM.INV(0.85, B1, B2)"

# B4:B68: normal density for each A-column x value, feeding the chart's
# y-values. B4 is entered on its own, then B5:B68 filled as one shared
# formula (matches how Excel groups a fill-down of the same formula).
$ws.Range("B4").Formula = "=NORM.DIST(A4,`$B`$1,`$B`$2,FALSE)"
$ws.Range("B5:B68").Formula = "=NORM.DIST(A5,`$B`$1,`$B`$2,FALSE)"

# Leave the selection where the author apparently ended up after entering H2
[void]$ws.Range("H3").Select()
